$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C7").Value = "MPZ2012S221AT000"
$ws.Range("G7").Value = "MPZ2012S221AT000"
$ws.Range("H7").Value = "TDK"
$ws.Range("I7").Value = "C76818"
